$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New text for F5 (EDU_LEVEL input_variables) - spacing around "occup_edu_h" tweaked
$newF5 = @'
school;occup_edu;occup_edu_2; occup_edu_h; educ; educ3
'@

# New text for H5 (EDU_LEVEL algorithm) - pipe separators tightened up (no surrounding spaces)
$newH5 = @'
case_when(
school %in% c(5,6)|occup_edu%in% c(7)|occup_edu_2  %in% c(7)|occup_edu_h%in% c(7)|educ %in% c(7)~ 7,                                                                                                                                               school%in% c(5,6)|occup_edu %in% c(6)|occup_edu_2  %in% c(6)|occup_edu_h%in% c(6)|educ%in% c(6)~ 6L,                                                                                                                                       school %in% c(5,6)|occup_edu%in% c(5)|occup_edu_2  %in% c(5)|occup_edu_h%in% c(5)|educ%in% c(5)~ 5L,                                                                                                                                 school %in% c(5,6)|occup_edu%in% c(4)|occup_edu_2 %in% c(4)|occup_edu_h%in% c(4)|educ%in% c(4)~ 4L,                                                                                                                         school %in% c(3,4)|occup_edu  %in% c(3)|occup_edu_2%in%c(3)|occup_edu_h%in% c(3)|educ %in% c(3)|educ3%in% c(3)~ 3L,                                                                                    school  %in% c(2)|occup_edu %in% c(2)|occup_edu_2 %in% c(2)|occup_edu_h%in% c(2)|educ %in% c(2)|educ3%in% c(2)~ 2L,                                                                                    school  %in% c(1)|occup_edu %in% c(1)|occup_edu_2 %in% c(1)|occup_edu_h%in% c(1)|educ %in% c(1)|educ3%in% c(1)~ 1L,                                       
school  %in% c(7)|occup_edu%in% c(8)|occup_edu_2  %in% c(8)|occup_edu_h%in% c(8)~ 9L,                                                                                                                                                                 TRUE ~ NA_integer_)
'@

$ws.Range("F5").Value = $newF5
$ws.Range("H5").Value = $newH5

# F18 keeps the same displayed text ("children"); Excel's shared-string table just
# gets reshuffled by the edit above, which this assignment reflects.
$ws.Range("F18").Value = "children"

# Column H got a bit wider now that the algorithm text changed, and is no longer
# flagged as "best fit" (a manual width was set instead).
$ws.Columns.Item(8).ColumnWidth = 78

# Leave the selection on the cell that was actually edited, scrolled back to the
# top of the sheet (matches the refreshed view after the fix).
$ws.Range("H5").Select()
